# Correcting market share tab for updated scenario 3s
$wb = $excel.ActiveWorkbook

$wsMarket = $wb.Worksheets.Item("MarketShare")

# --- Data changes on the MarketShare sheet ---
# Row 2 = "New Product A" market share: extend 1s into columns L:Y (2026-2039)
$wsMarket.Range("L2:Y2").Value = 1

# Row 3 = "Old Product B (SOC)" market share: it now only applies 2018-2025 (D:K);
# clear the trailing years L:Z (2026-2040) that used to hold 1s
$wsMarket.Range("L3:Z3").ClearContents()

# --- View / selection changes ---
# MarketShare becomes the active (selected) sheet/tab instead of Platform Coverage
$wsMarket.Select()
$wsMarket.Range("L4").Select()
$excel.ActiveWindow.ScrollColumn = 4
